# Update the "想去人数" (want-to-go count) figures in column F across all
# four sheets to reflect the latest scrape (gh-pages output regenerated at
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 263
$ws.Range("F4").Value  = 860
$ws.Range("F6").Value  = 430
$ws.Range("F7").Value  = 633
$ws.Range("F11").Value = 168
$ws.Range("F12").Value = 738
$ws.Range("F13").Value = 99
$ws.Range("F14").Value = 1868
$ws.Range("F15").Value = 390
$ws.Range("F16").Value = 4523
$ws.Range("F17").Value = 394

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 226
$ws.Range("F4").Value  = 47
$ws.Range("F6").Value  = 121
$ws.Range("F7").Value  = 490
$ws.Range("F13").Value = 101
$ws.Range("F14").Value = 43
$ws.Range("F15").Value = 6
$ws.Range("F19").Value = 25

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5396
$ws.Range("F3").Value = 344
$ws.Range("F4").Value = 312

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5396
$ws.Range("F4").Value  = 344
$ws.Range("F6").Value  = 312
$ws.Range("F7").Value  = 263
$ws.Range("F8").Value  = 226
$ws.Range("F9").Value  = 47
$ws.Range("F11").Value = 121
$ws.Range("F12").Value = 490
$ws.Range("F13").Value = 860
$ws.Range("F17").Value = 430
$ws.Range("F18").Value = 633
$ws.Range("F23").Value = 168
$ws.Range("F26").Value = 738
$ws.Range("F27").Value = 99
$ws.Range("F28").Value = 101
$ws.Range("F29").Value = 1868
$ws.Range("F30").Value = 390
$ws.Range("F31").Value = 4523
$ws.Range("F32").Value = 43
$ws.Range("F33").Value = 394
$ws.Range("F37").Value = 6
$ws.Range("F43").Value = 25

$wb.Save()
